# Append/update timestamp in the "ランサーズ" sheet: change A2:A12 from
# "2025-11-08 12:33:05" to "2025-11-08 12:41:49" (per commit message:
# "Append: 2025-11-08 12:41 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-08 12:41:49"

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
